$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 84
$ws1.Range("F5").Value = 19
$ws1.Range("F7").Value = 1692
$ws1.Range("F8").Value = 25
$ws1.Range("F11").Value = 1685
$ws1.Range("F13").Value = 85
$ws1.Range("F14").Value = 412
$ws1.Range("F21").Value = 274
$ws1.Range("F25").Value = 249

# Sheet "全部类型" (All types) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 84
$ws4.Range("F5").Value = 19
$ws4.Range("F7").Value = 1692
$ws4.Range("F9").Value = 25
$ws4.Range("F12").Value = 1685
$ws4.Range("F14").Value = 85
$ws4.Range("F15").Value = 412
$ws4.Range("F22").Value = 274
$ws4.Range("F26").Value = 249

$wb.Save()
